$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45880
$ws.Range("B2").Value = 109.7
$ws.Range("C2").Value = 100.43
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 99.89
$ws.Range("F2").Value = 99.18000000000001
$ws.Range("G2").Value = 99.89
$ws.Range("H2").Value = 99.65000000000001
$ws.Range("I2").Value = 105.43
$ws.Range("J2").Value = 100.61
$ws.Range("K2").Value = 63.6
$ws.Range("L2").Value = 45
$ws.Range("M2").Value = 22.51
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 15
$ws.Range("P2").Value = 10
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 44.5
$ws.Range("S2").Value = 57.97
$ws.Range("T2").Value = 90
$ws.Range("U2").Value = 111.24
$ws.Range("V2").Value = 135.09
$ws.Range("W2").Value = 144.41
$ws.Range("X2").Value = 133.51
$ws.Range("Y2").Value = 110.91
$ws.Range("Z2").Value = 80.09999999999999
$ws.Range("AB2").Value = 130.98
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 139.75
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 122.21
$ws.Range("AG2").Value = "9h-17h"
